$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (Volume number, week range)
$ws.Range("A8").Value = "Volume 30   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/17/2023  Through  4/23/2023"

# Plain numeric value updates
$ws.Range("L14").Value = -50
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 4
$ws.Range("K15").Value = 100
$ws.Range("M15").Value = -42.857142857142
$ws.Range("N15").Value = -42.857142857142
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 17
$ws.Range("H16").Value = -41.379310344827
$ws.Range("I16").Value = 74
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = -14.942528735632
$ws.Range("L16").Value = 85
$ws.Range("M16").Value = 8.823529411764
$ws.Range("N16").Value = -67.965367965368
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = 12.5
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = -27.272727272727
$ws.Range("I17").Value = 106
$ws.Range("J17").Value = 129
$ws.Range("K17").Value = -17.829457364341
$ws.Range("L17").Value = 26.190476190476
$ws.Range("M17").Value = 49.295774647887
$ws.Range("N17").Value = -39.772727272727
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = -35.294117647058
$ws.Range("L18").Value = 94.117647058823
$ws.Range("M18").Value = 3.125
$ws.Range("N18").Value = -77.397260273972
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 21.739130434782
$ws.Range("I19").Value = 119
$ws.Range("J19").Value = 133
$ws.Range("K19").Value = -10.526315789473
$ws.Range("L19").Value = 41.666666666666
$ws.Range("M19").Value = 105.172413793103
$ws.Range("N19").Value = 41.666666666666
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 20
$ws.Range("I20").Value = 17
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = -37.037037037037
$ws.Range("L20").Value = 6.25
$ws.Range("M20").Value = 41.666666666666
$ws.Range("N20").Value = -81.318681318681
$ws.Range("C21").Value = 22
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -14.705882352941
$ws.Range("I21").Value = 355
$ws.Range("J21").Value = 432
$ws.Range("K21").Value = -17.824074074074
$ws.Range("L21").Value = 42.570281124498
$ws.Range("M21").Value = 41.434262948207
$ws.Range("N21").Value = -52.156334231805
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 10
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 45
$ws.Range("K23").Value = 4.444444444444
$ws.Range("L23").Value = -9.615384615384
$ws.Range("M23").Value = 95.833333333333
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 13.636363636363
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 96
$ws.Range("H24").Value = -11.458333333333
$ws.Range("I24").Value = 392
$ws.Range("J24").Value = 333
$ws.Range("K24").Value = 17.717717717717
$ws.Range("L24").Value = 108.510638297872
$ws.Range("M24").Value = 30.666666666666
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 36
$ws.Range("G25").Value = 44
$ws.Range("H25").Value = -18.181818181818
$ws.Range("I25").Value = 150
$ws.Range("J25").Value = 173
$ws.Range("K25").Value = -13.294797687861
$ws.Range("L25").Value = 17.1875
$ws.Range("M25").Value = 20
$ws.Range("D26").Value = 1
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -33.333333333333
$ws.Range("I26").Value = 9
$ws.Range("J26").Value = 6
$ws.Range("K26").Value = 50
$ws.Range("L26").Value = 50
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 6
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 66.666666666666
$ws.Range("J28").Value = 12
$ws.Range("K28").Value = 8.333333333333
$ws.Range("M28").Value = 44.444444444444
$ws.Range("E29").Value = -100
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = 66.666666666666
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = 9.090909090909
$ws.Range("M29").Value = 50

# Cells changing from numeric style to text style "0" (shared string "0", style of C15)
$ws.Range("C14").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").Value = "'0"
$ws.Range("C15").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# Cells changing from numeric style to text style "***.*" (style of C15)
$ws.Range("E15").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("C15").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# L30: text style -> numeric style (matches K30 style)
$ws.Range("L30").Value = -100
$ws.Range("K30").Copy()
$ws.Range("L30").PasteSpecial(-4122)
